$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: text value changes from "3273112" to "3273114" (kept as text, not a number) ---
# Go through TEXT()+paste-values so the digit string stays a string (no number
# coercion) without touching the cell's number format / style.
$ws.Range("A1").Formula = '=TEXT(3273114,"0")'
$ws.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4163)

# --- B1: drop its numeric value (5); give it the new bold / wrap-text style ---
$ws.Range("B1").ClearContents()
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Size = 12
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").WrapText = $true

# --- New row 2, A2 holds the new part-number text "6SL32105BE211UV0" ---
$ws.Range("A2").Formula = '="6SL32105BE211UV0"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# --- Row heights: row 1 grows slightly, row 2 is tall enough for the wrapped text ---
$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 44

# --- Selection moves to A2 ---
[void]$ws.Range("A2").Select()
